# Natmi following Dr Hou advice
# Update Ligand-/Receptor-expressing cell counts (E, K: 1 -> 3) and the
# dependent expression / specificity statistics for all data rows (2-16)
# on the active (and only) worksheet of the NATMI LR-pairs export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "E" = 3; "G" = 1.139906333333333; "H" = 3.419719; "I" = 0.2178538649973528; "J" = 0.2178538649973527; "K" = 3; "M" = 8.591717333333333; "N" = 25.775152; "O" = 0.08461542565075156; "P" = 0.08461542565075157; "Q" = 9.793753002476445; "R" = 88.14377702228799; "S" = 0.01843379751641237; "T" = 0.01843379751641237 }
    3 = @{ "E" = 3; "G" = 1.139906333333333; "H" = 3.419719; "I" = 0.2178538649973528; "J" = 0.2178538649973527; "K" = 3; "M" = 16.543524; "N" = 49.630572; "O" = 0.1629286987355214; "P" = 0.1629286987355215; "Q" = 18.858067783252; "R" = 169.722610049268; "S" = 0.03549464673852264; "T" = 0.03549464673852264 }
    4 = @{ "E" = 3; "G" = 1.139906333333333; "H" = 3.419719; "I" = 0.2178538649973528; "J" = 0.2178538649973527; "K" = 3; "M" = 35.53801933333333; "N" = 106.614058; "O" = 0.3499957594051789; "P" = 0.3499957594051789; "Q" = 40.5100133121891; "R" = 364.5901198097019; "S" = 0.07624792891910179; "T" = 0.07624792891910179 }
    5 = @{ "E" = 3; "G" = 1.139906333333333; "H" = 3.419719; "I" = 0.2178538649973528; "J" = 0.2178538649973527; "K" = 3; "M" = 25.92369633333334; "N" = 77.771089; "O" = 0.2553092140468263; "P" = 0.2553092140468263; "Q" = 29.55058563377678; "R" = 265.955270703991; "S" = 0.05562009904953754; "T" = 0.05562009904953753 }
    6 = @{ "E" = 3; "G" = 1.139906333333333; "H" = 3.419719; "I" = 0.2178538649973528; "J" = 0.2178538649973527; "K" = 3; "M" = 14.94147133333333; "N" = 44.824414; "O" = 0.1471509021617218; "P" = 0.1471509021617218; "Q" = 17.03187780218511; "R" = 153.286900219666; "S" = 0.0320573927737784; "T" = 0.03205739277377841 }
    7 = @{ "E" = 3; "G" = 1.902924; "H" = 5.708772; "I" = 0.3636784322304457; "J" = 0.3636784322304456; "K" = 3; "M" = 8.591717333333333; "N" = 25.775152; "O" = 0.08461542565075156; "P" = 0.08461542565075157; "Q" = 16.349385114816; "R" = 147.144466033344; "S" = 0.03077280534317716; "T" = 0.03077280534317716 }
    8 = @{ "E" = 3; "G" = 1.902924; "H" = 5.708772; "I" = 0.3636784322304457; "J" = 0.3636784322304456; "K" = 3; "M" = 16.543524; "N" = 49.630572; "O" = 0.1629286987355214; "P" = 0.1629286987355215; "Q" = 31.481068864176; "R" = 283.329619777584; "S" = 0.05925365372148103; "T" = 0.05925365372148103 }
    9 = @{ "E" = 3; "G" = 1.902924; "H" = 5.708772; "I" = 0.3636784322304457; "J" = 0.3636784322304456; "K" = 3; "M" = 35.53801933333333; "N" = 106.614058; "O" = 0.3499957594051789; "P" = 0.3499957594051789; "Q" = 67.626149901864; "R" = 608.635349116776; "S" = 0.1272859090677797; "T" = 0.1272859090677797 }
    10 = @{ "E" = 3; "G" = 1.902924; "H" = 5.708772; "I" = 0.3636784322304457; "J" = 0.3636784322304456; "K" = 3; "M" = 25.92369633333334; "N" = 77.771089; "O" = 0.2553092140468263; "P" = 0.2553092140468263; "Q" = 49.330823921412; "R" = 443.977415292708; "S" = 0.09285045469853707; "T" = 0.09285045469853706 }
    11 = @{ "E" = 3; "G" = 1.902924; "H" = 5.708772; "I" = 0.3636784322304457; "J" = 0.3636784322304456; "K" = 3; "M" = 14.94147133333333; "N" = 44.824414; "O" = 0.1471509021617218; "P" = 0.1471509021617218; "Q" = 28.432484395512; "R" = 255.892359559608; "S" = 0.05351560939947067; "T" = 0.05351560939947068 }
    12 = @{ "E" = 3; "G" = 2.189605333333333; "H" = 6.568816; "I" = 0.4184677027722017; "J" = 0.4184677027722016; "K" = 3; "M" = 8.591717333333333; "N" = 25.775152; "O" = 0.08461542565075156; "P" = 0.08461542565075157; "Q" = 18.81247009555911; "R" = 169.312230860032; "S" = 0.03540882279116203; "T" = 0.03540882279116203 }
    13 = @{ "E" = 3; "G" = 2.189605333333333; "H" = 6.568816; "I" = 0.4184677027722017; "J" = 0.4184677027722016; "K" = 3; "M" = 16.543524; "N" = 49.630572; "O" = 0.1629286987355214; "P" = 0.1629286987355215; "Q" = 36.223788382528; "R" = 326.014095442752; "S" = 0.06818039827551778; "T" = 0.06818039827551778 }
    14 = @{ "E" = 3; "G" = 2.189605333333333; "H" = 6.568816; "I" = 0.4184677027722017; "J" = 0.4184677027722016; "K" = 3; "M" = 35.53801933333333; "N" = 106.614058; "O" = 0.3499957594051789; "P" = 0.3499957594051789; "Q" = 77.81423666836977; "R" = 700.328130015328; "S" = 0.1464619214182974; "T" = 0.1464619214182974 }
    15 = @{ "E" = 3; "G" = 2.189605333333333; "H" = 6.568816; "I" = 0.4184677027722017; "J" = 0.4184677027722016; "K" = 3; "M" = 25.92369633333334; "N" = 77.771089; "O" = 0.2553092140468263; "P" = 0.2553092140468263; "Q" = 56.76266375118045; "R" = 510.863973760624; "S" = 0.1068386602987517; "T" = 0.1068386602987517 }
    16 = @{ "E" = 3; "G" = 2.189605333333333; "H" = 6.568816; "I" = 0.4184677027722017; "J" = 0.4184677027722016; "K" = 3; "M" = 14.94147133333333; "N" = 44.824414; "O" = 0.1471509021617218; "P" = 0.1471509021617218; "Q" = 32.71592531931378; "R" = 294.443327873824; "S" = 0.06157789998847272; "T" = 0.06157789998847273 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
